# Apply the "tiny toss" ability example changes to the 技能表 workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before AC (column 29). This shifts the existing
#    AC:AJ columns (damage type thru precache) one slot to the right,
#    becoming AD:AK, and opens up a blank AC column for the new
#    "AbilityUnitDamageType" field.
# ---------------------------------------------------------------------------
$ws.Columns("AC:AC").Insert(-4161) | Out-Null

# ---------------------------------------------------------------------------
# 2. Populate the new AC column (rows 1-4) with the ability-damage-type
#    key/value/description triple.
# ---------------------------------------------------------------------------
$ws.Range("AC1").Value = "技能伤害类型"
$ws.Range("AC2").Value = "AbilityUnitDamageType"
$ws.Range("AC3").Value = "DAMAGE_TYPE_MAGICAL"
$ws.Range("AC4").Value = "DAMAGE_TYPE_PURE"

# ---------------------------------------------------------------------------
# 3. Widen the (now shifted) last column AK, which holds the "预载资源"/
#    Precache block, to fit its multi-line text.
# ---------------------------------------------------------------------------
$ws.Columns("AK:AK").ColumnWidth = (16.640625 - 5/7)

# ---------------------------------------------------------------------------
# 4. Add the new "tiny_toss_x" example ability as row 5.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "tiny_toss_x"
$ws.Range("B5").Value = "投掷"
$ws.Range("C5").Value = "投掷的技能描述！（懒得去找了）"
$ws.Range("D5").Value = "ability_lua"
$ws.Range("E5").Formula = '="examples/abilities/"&A5&".lua"'

$ws.Range("G5").Value = "duration 1.1"
$ws.Range("H5").Value = "持续事件"
$ws.Range("I5").Value = "grab_radius 275"

$ws.Range("K5").Value = "radius 275"
$ws.Range("L5").Value = "范围"
$ws.Range("M5").Value = "bonus_damage_pct 30 60 90"

$ws.Range("N5").Value = "被投掷者额外伤害"
$ws.Range("O5").Value = "toss_damage 200 250 300 350"
$ws.Range("P5").Value = "投掷伤害"

$ws.Range("X5").Value = 4
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 70

$ws.Range("AA5").Value = "tiny_toss"
$ws.Range("AC5").Value = "DAMAGE_TYPE_MAGICAL"

$ws.Range("AI5").Value = 1200
$ws.Range("AJ5").Value = 0

$ws.Range("AK5").Value = "{`n            ""particle""  ""particles/units/heroes/hero_tiny/tiny_toss_blur.vpcf""`n`t`t}"

# ---------------------------------------------------------------------------
# 5. The two "precache block" cells (AK4 old-axe, AK5 new-tiny) wrap their
#    long text onto multiple lines.
# ---------------------------------------------------------------------------
$ws.Range("AK4").WrapText = $true
$ws.Range("AK5").WrapText = $true

# Keep row 5's height in line with the other data rows (multi-line content
# would otherwise auto-expand it).
$ws.Rows("5:5").RowHeight = 14.25

# ---------------------------------------------------------------------------
# 6. Re-create the sheet view: freeze the header rows/first column, and
#    leave the selection sitting on M5 (the new damage-percent cell).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("B3").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("M5").Select() | Out-Null
